$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-24 Wednesday" "2024-01-25 Thursday"

Replace-Text "767÷5=" "726÷5="
Replace-Text "998÷2=" "839÷3="
Replace-Text "594÷9=" "903÷3="
Replace-Text "347÷3=" "735÷5="
Replace-Text "569÷9=" "716÷5="
Replace-Text "766÷4=" "430÷2="
Replace-Text "895÷8=" "408÷4="
Replace-Text "888÷3=" "695÷5="
Replace-Text "623÷3=" "982÷4="
Replace-Text "776÷5=" "540÷2="
Replace-Text "774÷3=" "980÷9="
Replace-Text "841÷6=" "357÷4="
Replace-Text "163÷6=" "620÷7="
Replace-Text "370÷9=" "410÷5="
Replace-Text "945÷4=" "968÷9="
Replace-Text "367÷2=" "974÷7="
Replace-Text "838÷5=" "857÷7="
Replace-Text "773÷3=" "568÷2="
Replace-Text "110÷9=" "515÷9="
Replace-Text "197÷8=" "723÷5="
Replace-Text "429÷6=" "908÷3="
Replace-Text "337÷4=" "130÷3="
Replace-Text "394÷7=" "733÷2="
Replace-Text "208÷8=" "877÷6="
Replace-Text "126÷5=" "299÷5="
